# Insert a new column before column A to hold the source file name,
# shifting all existing data (previously in columns A:J) one column to
# the right (now in columns B:K).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

# Give the new header cell the same formatting as the rest of the
# header row (bold, centered, bordered) before filling in its value.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Header for the new column.
$ws.Range("A1").Value = "Source File"

# The data rows came from two source CSV files. Only the first row of
# each file's block carries the file name; the remaining rows in that
# block are left blank, matching the original combination behavior.
$ws.Range("A2").Value = "SampleCSVFile_2kb_2.csv"
$ws.Range("A3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = ""
$ws.Range("A8").Value = ""
$ws.Range("A9").Value = ""

$ws.Range("A10").Value = "SampleCSVFile_2kb.csv"
$ws.Range("A11").Value = ""
$ws.Range("A12").Value = ""
$ws.Range("A13").Value = ""
$ws.Range("A14").Value = ""
$ws.Range("A15").Value = ""
$ws.Range("A16").Value = ""
$ws.Range("A17").Value = ""
